$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1100, 12.03800201416016, 1.239210605621338, 7.101234912872314, 18.45488166809082, 18476, "11"),
    @(1200, 12.04506397247314, 1.251126050949097, 6.868865013122559, 19.94649124145508, 18568, "11"),
    @(1300, 12.01037406921387, 1.252374768257141, 6.766349315643311, 20.14126968383789, 18438, "11"),
    @(1400, 11.99558448791504, 1.263855814933777, 7.012387275695801, 19.18616104125977, 18366, "11"),
    @(1500, 11.98649787902832, 1.29723048210144,  6.56302547454834,  20.06096649169922, 18392, "11")
)

$startRow = 21
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    # Column G ("Month") is stored as text ("11") in the source data, like
    # the existing rows above. Copying an existing text cell (G20) preserves
    # the text cell type instead of Excel auto-coercing "11" to a number.
    $ws.Range("G20").Copy($ws.Cells.Item($row, 7))
}
